$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.202.42'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '3.513.64'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '594.85'
$ws.Range("E5").Value = '  +0.26%  '
$ws.Range("D6").Value = '172.68'
$ws.Range("E6").Value = '  +1.44%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +2.20%  '
$ws.Range("E9").Value = '  +6.70%  '
$ws.Range("E10").Value = '  +0.35%  '
$ws.Range("D11").Value = '0.436'
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("D12").Value = '4.128.62'
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("D14").Value = '28.69'
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("D16").Value = '67.183.05'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").Value = '3.515.56'
$ws.Range("E17").Value = '  -0.43%  '
$ws.Range("D18").Value = '6.33'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").Value = '396.28'
$ws.Range("E20").Value = '  -0.46%  '
$ws.Range("E21").Value = '  +0.56%  '
$ws.Range("D22").Value = '73.39'
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("E24").Value = '  +1.44%  '
$ws.Range("E25").Value = '  -4.69%  '
$ws.Range("D26").Value = '10.29'
$ws.Range("E26").Value = '  +1.86%  '
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("E29").Value = '  -2.17%  '
$ws.Range("E30").Value = '  -1.58%  '
$ws.Range("E31").Value = '  +0.84%  '
$ws.Range("E32").Value = '  +1.83%  '
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("E34").Value = '  +2.21%  '
$ws.Range("D35").Value = '163.87'
$ws.Range("E35").Value = '  +1.14%  '
$ws.Range("D36").Value = '0.894'
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("D37").Value = '1.92'
$ws.Range("E37").Value = '  -1.41%  '
$ws.Range("D38").Value = '6.90'
$ws.Range("E38").Value = '  +2.71%  '
$ws.Range("D39").Value = '4.73'
$ws.Range("E39").Value = '  +0.63%  '
$ws.Range("E40").Value = '  -0.57%  '
$ws.Range("D41").Value = '27.37'
$ws.Range("E41").Value = '  +1.20%  '
$ws.Range("D42").Value = '26.39'
$ws.Range("E42").Value = '  -0.20%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.813.30'
$ws.Range("E43").Value = '  -0.58%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = '2.62'
$ws.Range("E44").Value = '  +1.84%  '
$ws.Range("D45").Value = '42.91'
$ws.Range("E45").Value = '  -1.49%  '
$ws.Range("E46").Value = '  -2.26%  '
$ws.Range("D47").Value = '342.15'
$ws.Range("E47").Value = '  -3.09%  '
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("D49").Value = '33.98'
$ws.Range("E49").Value = '  +1.92%  '
$ws.Range("D50").Value = '6.52'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("E51").Value = '  -0.33%  '
